# Hjemme passive updated meanEMG legmaxROM
# Updates the B:E columns of rows 1-3 (header "15/16" groupings + the two
# recomputed subject rows) and moves the sheet selection to the new B1:E3
# data block, matching the latest export of this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 (header) values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values
$ws.Range("B2").Value = 36.719898939644025
$ws.Range("C2").Value = 7.9019510033225515
$ws.Range("D2").Value = 11.458541932058267
$ws.Range("E2").Value = 1.983197291008338

# Row 3 values
$ws.Range("B3").Value = 52.404265362008033
$ws.Range("C3").Value = 6.8979008167237339
$ws.Range("D3").Value = -10.319453557895372
$ws.Range("E3").Value = 12.786249941062096

# The workbook's saved selection now spans just the B1:E3 block
$ws.Range("B1:E3").Select()
